$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2 had a standalone formula (not part of the shared group) -> update alone.
$ws.Range("B2").Formula = "=A2^(1/2)"

# B3:B22 was a shared-formula group (master in B3, ref B3:B22, si=0).
# Re-applying the new formula to the same range keeps it as one shared
# group with the master formula text updated in place.
$ws.Range("B3:B22").Formula = "=A3^(1/2)"

# Update the selection to match the saved state (B2:B22 selected, active cell B2).
$ws.Range("B2:B22").Select() | Out-Null
